$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14/15 swap: WrappedliquidstakedEther2.0 <-> Avalanche ---
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'32.16"
$ws.Range("E14").Value = "  +0.08%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.115.91"
$ws.Range("E15").Value = "  -0.22%  "

# --- Row 31/32 swap: Binance-PegBSC-USD <-> Fetch.AI ---
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.55"
$ws.Range("E31").Value = "  -6.74%  "

$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.11%  "

# --- Remaining D/E value updates ---
$ws.Range("D2").Value = "67.434.27"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "3.524.42"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'613.95"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'151.66"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("D7").Value = "3.524.16"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.479"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D11").Value = "'7.12"
$ws.Range("E11").Value = "  +3.08%  "
$ws.Range("D12").Value = "'0.426"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "'0.0000221"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D16").Value = "3.514.96"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "67.395.32"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "'6.40"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").Value = "'15.26"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "'445.20"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").Value = "'9.48"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").Value = "'0.626"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("D24").Value = "'77.39"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").Value = "'0.0000131"
$ws.Range("E25").Value = "  +10.38%  "
$ws.Range("D26").Value = "3.664.16"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "'10.28"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").Value = "'8.39"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").Value = "'2.51"
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("D33").Value = "'0.165"
$ws.Range("E33").Value = "  +4.45%  "
$ws.Range("D34").Value = "'25.91"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'6.18"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").Value = "3.511.20"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "'1.85"
$ws.Range("E37").Value = "  -2.86%  "
$ws.Range("D38").Value = "'8.02"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "'177.62"
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  +3.77%  "
$ws.Range("D43").Value = "'0.0883"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").Value = "'5.45"
$ws.Range("E44").Value = "  -2.69%  "
$ws.Range("D45").Value = "'0.882"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").Value = "'28.43"
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("D47").Value = "'45.12"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("D48").Value = "'2.63"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("E49").Value = "  +3.68%  "
$ws.Range("D50").Value = "'7.59"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -1.46%  "
